$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new blank paragraph right before the paragraph that starts
#    "EL PAGO DEBERÁ REALIZARSE DE LUNES A SÁBADO..." (same list/para
#    formatting is inherited automatically from that paragraph).
# ---------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("EL PAGO DEBERÁ REALIZARSE DE LUNES A SÁBADO", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $para1 = $r1.Paragraphs(1).Range
    $para1.InsertParagraphBefore()
}

# ---------------------------------------------------------------------
# 2) Split "ACUDIR A TODAS LAS DILIGENCIAS NECESARIAS PARA LA CONCLUSIÓN
#    DE LA NEGOCIACIÓN." into two runs right before "CONCLUSIÓN" (same
#    text / visual formatting, just two separate runs).
# ---------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("CONCLUSIÓN DE LA NEGOCIACIÓN.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $sub2 = $d.Range($r2.Start, $r2.End)
    $sub2.Font.Name = "Arial Narrow"
    $sub2.Font.Color = 0
}

# ---------------------------------------------------------------------
# 3) Merge the run-fragmented sentence about the 20% conventional
#    penalty back into a single run and drop the stray white highlight
#    that was applied to the middle fragment.
# ---------------------------------------------------------------------
$r3 = $d.Content
$oldText3 = "IMPUTABLE POR CUALQUIERA DE ELLAS, LA RESPONSABLE DEBERÁ CUBRIR COMO PENA CONVENCIONAL A LA OTRA, LA CANTIDAD EQUIVALENTE"
$found3 = $r3.Find.Execute($oldText3, $true, $false, $false, $false, $false, $true, 1, $false, $oldText3, 2)

# ---------------------------------------------------------------------
# 4) " QUEDARÁ OBLIGADO A CUBRIR LA PENA CONVENCIONAL" -> insert the
#    {{SEXO_5}} gender placeholder between "OBLIGAD" and the "O" so the
#    ending becomes templated, matching the {{SEXO_x}} convention used
#    elsewhere in the document.
# ---------------------------------------------------------------------
$r4 = $d.Content
$oldText4 = " QUEDARÁ OBLIGADO A CUBRIR LA PENA CONVENCIONAL"
$newText4 = " QUEDARÁ OBLIGAD{{SEXO_5}} A CUBRIR LA PENA CONVENCIONAL"
$found4 = $r4.Find.Execute($oldText4, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $start4 = $r4.Start
    $r4.Text = $newText4

    $prefixLen = ([string]" QUEDARÁ OBLIGAD").Length
    $tokenLen = ([string]"{{SEXO_5}}").Length
    $tokenStart = $start4 + $prefixLen
    $tokenEnd = $tokenStart + $tokenLen
    $subToken = $d.Range($tokenStart, $tokenEnd)
    $subToken.Font.Name = "Arial Narrow"
    $subToken.Font.Color = 0
}
